# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) across the per-job sheets with freshly pulled
# Universalis price data. Values only - no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1027.5
$ws.Range("I31").Value = 555
$ws.Range("K31").Value = 1665
$ws.Range("M31").Value = -1435

$ws.Range("H98").Value = 802.65
$ws.Range("I98").Value = 578
$ws.Range("J98").Value = 1701.25
$ws.Range("K98").Value = 578
$ws.Range("L98").Value = 1701.25
$ws.Range("M98").Value = 920
$ws.Range("N98").Value = -4697.25

$ws.Range("H122").Value = 802.65
$ws.Range("I122").Value = 578
$ws.Range("J122").Value = 1701.25
$ws.Range("K122").Value = 1734
$ws.Range("L122").Value = 5103.75
$ws.Range("M122").Value = 716
$ws.Range("N122").Value = -10003.75

$ws.Range("H132").Value = 4613.8423
$ws.Range("I132").Value = 4291.125
$ws.Range("J132").Value = 6335
$ws.Range("K132").Value = 12873.375
$ws.Range("L132").Value = 19005
$ws.Range("M132").Value = -10343.375
$ws.Range("N132").Value = -24065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 448668.44
$ws.Range("I32").Value = 477324.06
$ws.Range("J32").Value = 28385.666
$ws.Range("K32").Value = 477324.06
$ws.Range("L32").Value = 28385.666
$ws.Range("M32").Value = -477037.06
$ws.Range("N32").Value = -28959.666

$ws.Range("H45").Value = 3922.2727
$ws.Range("I45").Value = 3842.2
$ws.Range("J45").Value = 3989
$ws.Range("K45").Value = 3842.2
$ws.Range("L45").Value = 3989
$ws.Range("M45").Value = -3465.2
$ws.Range("N45").Value = -4743

$ws.Range("H61").Value = 16395390
$ws.Range("I61").Value = 21740390
$ws.Range("K61").Value = 21740390
$ws.Range("M61").Value = -21740178

$ws.Range("H74").Value = 1139.5641
$ws.Range("I74").Value = 886.86365
$ws.Range("J74").Value = 1466.5883
$ws.Range("K74").Value = 886.86365
$ws.Range("L74").Value = 1466.5883
$ws.Range("M74").Value = -12.86365000000001
$ws.Range("N74").Value = -3214.5883

$ws.Range("H77").Value = 1139.5641
$ws.Range("I77").Value = 886.86365
$ws.Range("J77").Value = 1466.5883
$ws.Range("K77").Value = 4434.31825
$ws.Range("L77").Value = 7332.941499999999
$ws.Range("M77").Value = -66.31825000000026
$ws.Range("N77").Value = -16068.9415

$ws.Range("H110").Value = 54986.59
$ws.Range("I110").Value = 66190.07000000001
$ws.Range("K110").Value = 66190.07000000001
$ws.Range("M110").Value = -64145.07000000001

$ws.Range("H132").Value = 2266.1885
$ws.Range("I132").Value = 1403.0392
$ws.Range("J132").Value = 4711.778
$ws.Range("K132").Value = 4209.1176
$ws.Range("L132").Value = 14135.334
$ws.Range("M132").Value = -1679.1176
$ws.Range("N132").Value = -19195.334

$ws.Range("H136").Value = 16395390
$ws.Range("I136").Value = 21740390
$ws.Range("K136").Value = 65221170
$ws.Range("M136").Value = -65218620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10419153
$ws.Range("I105").Value = 10419153
$ws.Range("K105").Value = 10419153
$ws.Range("M105").Value = -10417406

$ws.Range("H107").Value = 1037.1818
$ws.Range("I107").Value = 845.5
$ws.Range("J107").Value = 1548.3334
$ws.Range("K107").Value = 845.5
$ws.Range("L107").Value = 1548.3334
$ws.Range("M107").Value = 1074.5
$ws.Range("N107").Value = -5388.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1717.1333
$ws.Range("I99").Value = 1289.2222
$ws.Range("J99").Value = 1900.5238
$ws.Range("K99").Value = 1289.2222
$ws.Range("L99").Value = 1900.5238
$ws.Range("M99").Value = 208.7778000000001
$ws.Range("N99").Value = -4896.5238

$ws.Range("H126").Value = 1717.1333
$ws.Range("I126").Value = 1289.2222
$ws.Range("J126").Value = 1900.5238
$ws.Range("K126").Value = 3867.6666
$ws.Range("L126").Value = 5701.5714
$ws.Range("M126").Value = -1397.6666
$ws.Range("N126").Value = -10641.5714

$ws.Range("H132").Value = 18818192
$ws.Range("I132").Value = 29412624
$ws.Range("J132").Value = 5953525
$ws.Range("K132").Value = 88237872
$ws.Range("L132").Value = 17860575
$ws.Range("M132").Value = -88235342
$ws.Range("N132").Value = -17865635

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 32000
$ws.Range("J34").Value = 32000
$ws.Range("L34").Value = 32000
$ws.Range("N34").Value = -32536

$ws.Range("H76").Value = 32000
$ws.Range("J76").Value = 32000
$ws.Range("L76").Value = 32000
$ws.Range("N76").Value = -32630

$ws.Range("H79").Value = 32000
$ws.Range("J79").Value = 32000
$ws.Range("L79").Value = 32000
$ws.Range("N79").Value = -34184

$ws.Range("H113").Value = 2523.6667
$ws.Range("I113").Value = 2537.3845
$ws.Range("J113").Value = 2488
$ws.Range("K113").Value = 2537.3845
$ws.Range("L113").Value = 2488
$ws.Range("M113").Value = -367.3845000000001
$ws.Range("N113").Value = -6828

$ws.Range("H132").Value = 3052.3264
$ws.Range("I132").Value = 2786.8572
$ws.Range("J132").Value = 3716
$ws.Range("K132").Value = 8360.571599999999
$ws.Range("L132").Value = 11148
$ws.Range("M132").Value = -5830.571599999999
$ws.Range("N132").Value = -16208

$ws.Range("H134").Value = 32986.445
$ws.Range("J134").Value = 32986.445
$ws.Range("L134").Value = 98959.33499999999
$ws.Range("N134").Value = -104029.335

$ws.Range("H136").Value = 32217.334
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 32217.334
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 96652.00199999999
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -101752.002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 1850
$ws.Range("K46").Value = 1850
$ws.Range("M46").Value = -1662

$ws.Range("H101").Value = 25313.125
$ws.Range("J101").Value = 25313.125
$ws.Range("L101").Value = 25313.125
$ws.Range("N101").Value = -31803.125

$ws.Range("H132").Value = 2338.1875
$ws.Range("I132").Value = 2137.3225
$ws.Range("J132").Value = 2704.4707
$ws.Range("K132").Value = 6411.967500000001
$ws.Range("L132").Value = 8113.4121
$ws.Range("M132").Value = -3881.967500000001
$ws.Range("N132").Value = -13173.4121

$ws.Range("H136").Value = 6174418.5
$ws.Range("I136").Value = 1380.8
$ws.Range("J136").Value = 23811668
$ws.Range("K136").Value = 4142.4
$ws.Range("L136").Value = 71435004
$ws.Range("M136").Value = -1592.4
$ws.Range("N136").Value = -71440104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 33000
$ws.Range("J75").Value = 33000
$ws.Range("L75").Value = 33000
$ws.Range("N75").Value = -34872

$ws.Range("H78").Value = 33000
$ws.Range("J78").Value = 33000
$ws.Range("L78").Value = 99000
$ws.Range("N78").Value = -108360

$ws.Range("H103").Value = 55200.668
$ws.Range("J103").Value = 55200.668
$ws.Range("L103").Value = 55200.668
$ws.Range("N103").Value = -57544.668

$ws.Range("H132").Value = 2875418.8
$ws.Range("I132").Value = 1782.225
$ws.Range("J132").Value = 9261278
$ws.Range("K132").Value = 5346.674999999999
$ws.Range("L132").Value = 27783834
$ws.Range("M132").Value = -2816.674999999999
$ws.Range("N132").Value = -27788894

$ws.Range("H136").Value = 2179.8135
$ws.Range("I136").Value = 1725.54
$ws.Range("K136").Value = 5176.62
$ws.Range("M136").Value = -2626.62

$ws.Range("H137").Value = 88007.5
$ws.Range("J137").Value = 88007.5
$ws.Range("L137").Value = 88007.5
$ws.Range("N137").Value = -98207.5
